# [Fix] item data classification fix: renumber a handful of item IDs
# (10202005-10202008 -> 10102005-10102008, 10302000/10303000/10304000/
#  10306000/10307000/10308000/10309000 -> 10132000/10133000/10134000/
#  10136000/10137000/10138000/10139000) across the ItemDatas, Stats and
# Recipe sheets, including the materials_string lists on Recipe that
# reference those ids.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "ItemDatas": id column (A) renumbering
# ---------------------------------------------------------------
$itemDatas = $wb.Worksheets.Item("ItemDatas")

$itemDatas.Range("A6").Value = 10102005
$itemDatas.Range("A7").Value = 10102006
$itemDatas.Range("A8").Value = 10102007
$itemDatas.Range("A9").Value = 10102008
$itemDatas.Range("A10").Value = 10132000
$itemDatas.Range("A11").Value = 10133000
$itemDatas.Range("A12").Value = 10134000
$itemDatas.Range("A13").Value = 10136000
$itemDatas.Range("A14").Value = 10137000
$itemDatas.Range("A15").Value = 10138000
$itemDatas.Range("A16").Value = 10139000

# ---------------------------------------------------------------
# Sheet "Stats": id column (A) renumbering (mirrors ItemDatas)
# ---------------------------------------------------------------
$stats = $wb.Worksheets.Item("Stats")

$stats.Range("A14").Value = 10102005
$stats.Range("A15").Value = 10102006
$stats.Range("A16").Value = 10102007
$stats.Range("A17").Value = 10102008
$stats.Range("A18").Value = 10132000
$stats.Range("A19").Value = 10133000
$stats.Range("A20").Value = 10134000
$stats.Range("A21").Value = 10136000
$stats.Range("A22").Value = 10137000
$stats.Range("A23").Value = 10138000
$stats.Range("A24").Value = 10139000

# ---------------------------------------------------------------
# Sheet "Recipe": craftingID (B) and materials_string (D) columns
# ---------------------------------------------------------------
$recipe = $wb.Worksheets.Item("Recipe")

$recipe.Range("B13").Value = 10132000
$recipe.Range("D13").Value = 10136000

$recipe.Range("B14").Value = 10133000
$recipe.Range("D14").Value = 10102005

$recipe.Range("B15").Value = 10134000
$recipe.Range("D15").Value = 10102007

$recipe.Range("B16").Value = 10136000
$recipe.Range("B17").Value = 10137000
$recipe.Range("B18").Value = 10138000
$recipe.Range("B19").Value = 10139000

# Materials-string lists (written in this order so the rebuilt shared
# string table lines up with the new entries appended by the authors'
# original edit).
$recipe.Range("D20").Value = "10102006, 10202002"
$recipe.Range("D19").Value = "10102006, 10102007"
$recipe.Range("D18").Value = "10102005, 10102007"
$recipe.Range("D16").Value = "10202006, 10102008"
$recipe.Range("D17").Value = "10139000, 10136000, 10102006"

# ---------------------------------------------------------------
# Restore / update each sheet's active selection
# ---------------------------------------------------------------
$itemDatas.Range("A6:A9").Select() | Out-Null
$recipe.Range("D17").Select() | Out-Null
$stats.Range("I20").Select() | Out-Null
